$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append below the existing header (row 1) and first
# data row (row 2): rows 3-5.
$rows = @(
    @("11/04/2022", "12:19", "12:26", 6, 5, 1, 1, 3, "EURUSD", 8.69, "WIN"),
    @("11/04/2022", "12:26", "12:30", 3, 0, 3, 3, 0, "EURUSD", -14, "LOSS"),
    @("11/04/2022", "12:31", "12:38", 6, 3, 3, 2, 2, "EURUSD", 2.699999999999999, "LOSS")
)

$startRow = 3
$endRow = $startRow + $rows.Length - 1

# Column A ("Data ") holds a dd/mm/yyyy-looking string that must be stored
# as literal text (like the existing rows), not auto-parsed into a date
# serial number. Temporarily force a text format on that range before
# typing the values, then clear the format override again afterwards so
# the cells end up with no explicit style - matching row 2's plain cells.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $rows[$i][0]
}
$dateRange.ClearFormats()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($col = 2; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
    }
}
